$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (44104) now holds what used to be row 3's hours + activity note
$ws.Range("B4").Value = 0.65844907407407405
$ws.Range("C4").Value = 0.86927083333333333
$ws.Range("D4").Value = 0.020833333333333332
$ws.Range("E4").Formula = "=C4-B4-D4"
$ws.Range("F4").Value = "* Added nodes`n* Added node serialization`n* Added Undo/Redo (hopefully) and file change tracking`n* Added search tree for nodes`n* Researched the internals of ShaderGraph to learn how a bunch of things are done there, then reverse-engineered some of them"
$ws.Range("F4").WrapText = $true

# Row 3 (44103) is cleared back out to zero / no note
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("F3").ClearContents()

# Row 5 (44105) gets the new day's hours + activity note
$ws.Range("B5").Value = 0.47126157407407404
$ws.Range("C5").Value = 0.83059027777777772
$ws.Range("F5").Value = "* Spent a lot of time fixing the Undo/Redo functionality`n* Spent a bunch more time fixing the Undo/Redo functionality because apparently it's so impossible for my tiny brain to understand that you need to save the Undo state of the object before modifying it, and not afterwards leading to the Undo doing absolutely nothing. :/`n* Researched how ShaderGraph implements some features and reversed engineered that into something that works with my architecture (there are quite a few similar things in the codebase but I never just copy pasted things)`n* Added node connections (edges) and proper serialization and file persistancy for them`n* Added the ability to drag an edge/connection from a port and drop it somewhere to create a node"
$ws.Range("F5").WrapText = $true

# Give every data row (1-58) an explicit 15pt row height (matches the saved file)
for ($r = 1; $r -le 58; $r++) {
    $ws.Rows.Item($r).RowHeight = 15
}

$ws.Range("F5").Select()
